$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace date (B2) and added-time (H2) with the values that used to
# live in row 4. Use T("...") + paste-as-values so the date-like text is
# committed as a literal shared string instead of being auto-converted to a
# date serial number (which would also pull in a new number-format style).
$ws.Range("B2").Formula = '=T("11/10/2024")'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("H2").Value = "15:11:07"

# Row 3: replace date (B3) and added-time (H3) with the values that used to
# live in row 5.
$ws.Range("B3").Formula = '=T("12/10/2024")'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("H3").Value = "11:11:56"

# Drop the now-duplicate rows 4 and 5 entirely.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
